$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the "Emergencias" / "Equipos y más" rows (row 4)
$ws.Range("O4").Value = "Emergencias y rescate"
$ws.Range("P4").Value = "Equipos y mas"

# Update the "Prueba" / "Desc Prueba" rows (row 6), plus the associated ID value
$ws.Range("N6").Value = 57
$ws.Range("O6").Value = "Ingles B2"
$ws.Range("P6").Value = "Equipo de grammary"
